$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 22 data (21st candidate: Eduardo Artés)
$ws.Range("A22:D22").HorizontalAlignment = -4131
$ws.Range("F22").HorizontalAlignment = -4131

$ws.Range("A22").Value = 21
$ws.Range("B22").Value = "Eduardo Artés"
$ws.Range("C22").Value = "Unión Patriótica"
$ws.Range("D22").Value = "UPA"
$ws.Range("E22").Value = "https://www.latercera.com/la-tercera-pm/noticia/artes-mulet-y-me-o-los-otros-candidatos-presidenciales-de-la-oposicion/RLHFQHI27BABNGODECRWPYCO3U/"
$ws.Range("F22").Value = "."
$ws.Range("G22").Value = 0

$ws.Range("B20").Select()
